$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.78
$ws.Range("H2").Value = 3.15
$ws.Range("I2").Value = 4.8
$ws.Range("P2").Value = 1.57
$ws.Range("Q2").Value = 2.12
$ws.Range("R2").Value = 2.22
$ws.Range("S2").Value = 1.52
$ws.Range("U2").Value = 6.9
$ws.Range("W2").Value = 14
$ws.Range("Z2").Value = 6.3
$ws.Range("AA2").Value = 6.5
$ws.Range("AE2").Value = 25
$ws.Range("AG2").Value = 100
$ws.Range("H3").Value = 3.25
$ws.Range("H4").Value = 3.3
$ws.Range("G5").Value = 4.75
$ws.Range("H5").Value = 3.2
$ws.Range("T5").Value = 8.5
$ws.Range("U5").Value = 21
$ws.Range("V5").Value = 17
$ws.Range("AA5").Value = 6.5
$ws.Range("AB5").Value = 23
$ws.Range("AD5").Value = 5
$ws.Range("AF5").Value = 9.5
$ws.Range("G6").Value = 2.3
$ws.Range("N6").Value = 1.93
$ws.Range("O6").Value = 1.93
$ws.Range("P6").Value = 1.4
$ws.Range("Q6").Value = 2.75
$ws.Range("R6").Value = 1.73
$ws.Range("S6").Value = 2
$ws.Range("T6").Value = 8.5
$ws.Range("AH6").Value = 23
$ws.Range("G7").Value = 2.45
$ws.Range("I7").Value = 2.75
$ws.Range("W7").Value = 23
$ws.Range("AF7").Value = 11
$ws.Range("I8").Value = 5.75
$ws.Range("L8").Value = 1.36
$ws.Range("M8").Value = 3.2
$ws.Range("N8").Value = 2.1
$ws.Range("O8").Value = 1.73
$ws.Range("P8").Value = 1.44
$ws.Range("Q8").Value = 2.63
$ws.Range("AB8").Value = 21
$ws.Range("AD8").Value = 12
$ws.Range("G10").Value = 2.4
$ws.Range("I10").Value = 2.88
$ws.Range("T10").Value = 7.5
$ws.Range("AA10").Value = 6.5
$ws.Range("AE10").Value = 13
$ws.Range("AG10").Value = 29
$ws.Range("AI10").Value = 34
$ws.Range("H11").Value = 3.9
$ws.Range("I11").Value = 1.75
$ws.Range("L11").Value = 1.22
$ws.Range("M11").Value = 4.33
$ws.Range("N11").Value = 1.73
$ws.Range("O11").Value = 2.1
$ws.Range("P11").Value = 1.33
$ws.Range("Q11").Value = 3.25
$ws.Range("R11").Value = 1.67
$ws.Range("S11").Value = 2.1
$ws.Range("U11").Value = 23
$ws.Range("Z11").Value = 13
$ws.Range("AD11").Value = 8.5
$ws.Range("G12").Value = 1.7
$ws.Range("H12").Value = 3.75
$ws.Range("I12").Value = 4.5
$ws.Range("L12").Value = 1.3
$ws.Range("M12").Value = 3.4
$ws.Range("N12").Value = 2
$ws.Range("O12").Value = 1.85
$ws.Range("P12").Value = 1.4
$ws.Range("Q12").Value = 2.75
$ws.Range("R12").Value = 1.95
$ws.Range("S12").Value = 1.8
$ws.Range("T12").Value = 6.5
$ws.Range("X12").Value = 15
$ws.Range("Y12").Value = 29
$ws.Range("Z12").Value = 10
$ws.Range("AA12").Value = 7
$ws.Range("AD12").Value = 12
$ws.Range("AE12").Value = 23
$ws.Range("AJ12").Value = 301
$ws.Range("G14").Value = 3.25
$ws.Range("H14").Value = 3.3
$ws.Range("I14").Value = 2.15
$ws.Range("P14").Value = 1.44
$ws.Range("Q14").Value = 2.63
$ws.Range("T14").Value = 9
$ws.Range("V14").Value = 12
$ws.Range("X14").Value = 29
$ws.Range("AE14").Value = 10
$ws.Range("AH14").Value = 19
$ws.Range("AI14").Value = 29
$ws.Range("G15").Value = 2.45
$ws.Range("I15").Value = 2.63
$ws.Range("V15").Value = 10
$ws.Range("X15").Value = 21
$ws.Range("Y15").Value = 29
$ws.Range("AB15").Value = 13
$ws.Range("AD15").Value = 9
$ws.Range("AG15").Value = 26
$ws.Range("R16").Value = 1.67
$ws.Range("S16").Value = 2.1
$ws.Range("H17").Value = 3.25
$ws.Range("L17").Value = 1.36
$ws.Range("M17").Value = 3
$ws.Range("R17").Value = 1.91
$ws.Range("S17").Value = 1.91
$ws.Range("Z17").Value = 9
$ws.Range("AD17").Value = 7
$ws.Range("AJ17").Value = 301
$ws.Range("G18").Value = 1.85
$ws.Range("H18").Value = 3.6
$ws.Range("I18").Value = 3.9
$ws.Range("L18").Value = 1.36
$ws.Range("M18").Value = 3
$ws.Range("N18").Value = 2.15
$ws.Range("O18").Value = 1.67
$ws.Range("AF18").Value = 13
$ws.Range("AJ18").Value = 351
$ws.Range("H20").Value = 7.7
$ws.Range("I20").Value = 29
$ws.Range("N20").Value = 1.31
$ws.Range("O20").Value = 3.15
$ws.Range("R20").Value = 2.8
$ws.Range("S20").Value = 1.39
$ws.Range("T20").Value = 8.75
$ws.Range("U20").Value = 5.6
$ws.Range("V20").Value = 11.75
$ws.Range("X20").Value = 10.5
$ws.Range("Y20").Value = 40
$ws.Range("Z20").Value = 19
$ws.Range("AA20").Value = 19.5
$ws.Range("AE20").Value = 450
$ws.Range("AH20").Value = 500
$ws.Range("O21").Value = 2.18
$ws.Range("Y21").Value = 35
$ws.Range("Z21").Value = 14
$ws.Range("AA21").Value = 7.2
$ws.Range("AD21").Value = 7.1
$ws.Range("AE21").Value = 6.8
$ws.Range("AG21").Value = 9
$ws.Range("AI21").Value = 16.5
$ws.Range("AJ21").Value = 250
$ws.Range("N23").Value = 1.85
$ws.Range("O23").Value = 2
$ws.Range("G27").Value = 1.09
$ws.Range("H27").Value = 7
$ws.Range("I27").Value = 20
$ws.Range("N27").Value = 1.27
$ws.Range("O27").Value = 3.45
$ws.Range("R27").Value = 2.15
$ws.Range("S27").Value = 1.62
$ws.Range("T27").Value = 9
$ws.Range("U27").Value = 6
$ws.Range("Y27").Value = 27
$ws.Range("Z27").Value = 21
$ws.Range("AA27").Value = 15.5
$ws.Range("AB27").Value = 30
$ws.Range("AC27").Value = 110
$ws.Range("AD27").Value = 65
$ws.Range("AF27").Value = 60
$ws.Range("AH27").Value = 300
$ws.Range("G29").Value = 1.19
$ws.Range("H29").Value = 5.8
$ws.Range("I29").Value = 10.5
$ws.Range("L29").Value = 1.11
$ws.Range("M29").Value = 6
$ws.Range("R29").Value = 1.99
$ws.Range("S29").Value = 1.76
$ws.Range("T29").Value = 8.5
$ws.Range("U29").Value = 6.3
$ws.Range("V29").Value = 8.25
$ws.Range("W29").Value = 6.5
$ws.Range("X29").Value = 8.5
$ws.Range("Y29").Value = 20
$ws.Range("Z29").Value = 19
$ws.Range("AA29").Value = 11
$ws.Range("AB29").Value = 19.5
$ws.Range("AC29").Value = 70
$ws.Range("AD29").Value = 29
$ws.Range("AE29").Value = 70
$ws.Range("AF29").Value = 28
$ws.Range("AG29").Value = 250
$ws.Range("AH29").Value = 100
$ws.Range("AI29").Value = 70
$ws.Range("AJ29").Value = 400
$ws.Range("I31").Value = 4.75
$ws.Range("K31").Value = 13
$ws.Range("N31").Value = 1.8
$ws.Range("O31").Value = 2
$ws.Range("U31").Value = 8
$ws.Range("AA31").Value = 8
$ws.Range("AE31").Value = 26
$ws.Range("G32").Value = 1.45
$ws.Range("H32").Value = 3.9
$ws.Range("J32").Value = 1.07
$ws.Range("K32").Value = 8.5
$ws.Range("L32").Value = 1.33
$ws.Range("M32").Value = 3.25
$ws.Range("N32").Value = 2.05
$ws.Range("O32").Value = 1.75
$ws.Range("U32").Value = 6
$ws.Range("Z32").Value = 8.5
$ws.Range("AC32").Value = 81
$ws.Range("AD32").Value = 15
$ws.Range("AI32").Value = 67
$ws.Range("AJ32").Value = 900
$ws.Range("G33").Value = 2.75
$ws.Range("H33").Value = 2.85
$ws.Range("I33").Value = 2.75
$ws.Range("M33").Value = 2.77
$ws.Range("N33").Value = 2.15
$ws.Range("O33").Value = 1.62
$ws.Range("R33").Value = 1.78
$ws.Range("S33").Value = 1.93
$ws.Range("T33").Value = 7.6
$ws.Range("V33").Value = 10
$ws.Range("X33").Value = 25
$ws.Range("AA33").Value = 5.5
$ws.Range("AC33").Value = 65
$ws.Range("AD33").Value = 7.9
$ws.Range("AG33").Value = 35
$ws.Range("AI33").Value = 32
